$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "Case 3" section: a Heading1 paragraph containing the
# text "Case 3", immediately followed by the "バグ" / "なぜ、左側に
# 枝が伸びない？" paragraph, immediately followed by a centered
# paragraph that holds only an inline picture plus the (hidden)
# "_GoBack" bookmark.
# ------------------------------------------------------------------
$headingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text.TrimEnd() -eq "Case 3") {
        $headingIndex = $i
    }
}

if ($headingIndex -gt 0) {
    $headingPara = $d.Paragraphs.Item($headingIndex)
    $bugPara     = $d.Paragraphs.Item($headingIndex + 1)
    $picPara     = $d.Paragraphs.Item($headingIndex + 2)

    # Remove any inline picture living in the paragraph that follows
    # the "バグ" paragraph (keeps the bookmark that lives in the same
    # paragraph mark).
    $picRange = $picPara.Range
    for ($s = $d.InlineShapes.Count; $s -ge 1; $s--) {
        $shape = $d.InlineShapes.Item($s)
        if ($shape.Range.Start -ge $picRange.Start -and $shape.Range.End -le $picRange.End) {
            $shape.Delete()
        }
    }

    # Delete the "Case 3" heading paragraph and the "バグ" paragraph
    # entirely (including their paragraph marks), leaving the
    # (now-empty) picture paragraph as the final paragraph.
    $delRange = $d.Range($headingPara.Range.Start, $picPara.Range.Start)
    $delRange.Delete()

    # The surviving paragraph (formerly the picture paragraph) should
    # no longer be center-justified once the picture is gone.
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastPara.Alignment = 0
}
